# Apply changes described by the commit:
#  - Change existing "Yes" entries in A2 and A3 to "No"
#  - Add a new row (row 5) with values Yes / SMOKE, matching the style of row 4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 and A3 from "Yes" to "No"
$ws.Range("A2").Value = "No"
$ws.Range("A3").Value = "No"

# Add new child row 5: Yes / SMOKE, using the same style as row 4 (A4/B4)
$ws.Range("A5").Value = "Yes"
$ws.Range("B5").Value = "SMOKE"

$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122) | Out-Null

$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
